# Adds new export data (one more building, two more receipts, one more user)
# to the "users list export" workbook, matching the "add export excel in
# users list" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "اطلاعات کلی" (general info) - update the single data row (row 3)
# with the real building figures instead of the placeholder "test" data.
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item(1)

$wsInfo.Range("B3").Value = "ساختمان 1"
$wsInfo.Range("D3").Value = 2
$wsInfo.Range("F3").Value = 850120000
$wsInfo.Range("G3").Value = 3
$wsInfo.Range("I3").Value = "کرمان"

# ---------------------------------------------------------------------
# Sheet 2: "رسید ها" (receipts) - refresh existing row 3 with real data and
# append two new receipt rows (4 and 5), each with its own hyperlink to the
# receipt image.
# ---------------------------------------------------------------------
$wsReceipts = $wb.Worksheets.Item(2)

# Row 3: replace the placeholder test receipt with the first real receipt.
$wsReceipts.Range("A3").Value = "m3VNu2Aflb"
$wsReceipts.Range("B3").Value = "Fzl2 M2|9130009997"
$wsReceipts.Range("C3").Value = 120000
$wsReceipts.Range("D3").Value = 12
$wsReceipts.Range("E3").Value = 1
$wsReceipts.Range("F3").Value = "2023-10-18 13:16:08"
$wsReceipts.Range("G3").Value = "2023-09-24 00:00:00(101  روز پیش)"
$wsReceipts.Range("H3").Value = "awdwad"
$wsReceipts.Range("I3").Value = "awd"

$wsReceipts.Range("L3").Hyperlinks.Delete()
$wsReceipts.Range("L3").Value = "http://127.0.0.1:8000/media/images/users/9130009991/2023-10-18/KilPOgqb.png"
$wsReceipts.Hyperlinks.Add($wsReceipts.Range("L3"), "http://127.0.0.1:8000/media/images/users/9130009991/2023-10-18/KilPOgqb.png")
$wsReceipts.Range("L3").Style = "Hyperlink"

# Row 4: new receipt.
$wsReceipts.Range("A4").Value = "wDA2H4NjXQ"
$wsReceipts.Range("B4").Value = "test test|9130009999"
$wsReceipts.Range("C4").Value = 750000000
$wsReceipts.Range("D4").Value = 114750
$wsReceipts.Range("E4").Value = 1.7
$wsReceipts.Range("F4").Value = "2023-10-17 12:24:37"
$wsReceipts.Range("G4").Value = "2023-10-05 00:00:00(90  روز پیش)"
$wsReceipts.Range("H4").Value = "علی رضایی"
$wsReceipts.Range("I4").Value = "سامان"
$wsReceipts.Range("J4").Value = "WAWD"

$wsReceipts.Range("L4").Value = "http://127.0.0.1:8000/media/images/users/9130009999/2023-10-17/cohfPzI5.png"
$wsReceipts.Hyperlinks.Add($wsReceipts.Range("L4"), "http://127.0.0.1:8000/media/images/users/9130009999/2023-10-17/cohfPzI5.png")
$wsReceipts.Range("L4").Style = "Hyperlink"

# Row 5: new receipt.
$wsReceipts.Range("A5").Value = "yGTOyj4lrj"
$wsReceipts.Range("B5").Value = "test test|9130009999"
$wsReceipts.Range("C5").Value = 100000000
$wsReceipts.Range("D5").Value = 9200
$wsReceipts.Range("E5").Value = 1
$wsReceipts.Range("F5").Value = "2023-10-17 11:57:03"
$wsReceipts.Range("G5").Value = "2023-10-03 00:00:00(92  روز پیش)"
$wsReceipts.Range("H5").Value = "علی رضایی"
$wsReceipts.Range("I5").Value = "سامان"

$wsReceipts.Range("L5").Value = "http://127.0.0.1:8000/media/images/users/9130009999/2023-10-17/62MNmnro.png"
$wsReceipts.Hyperlinks.Add($wsReceipts.Range("L5"), "http://127.0.0.1:8000/media/images/users/9130009999/2023-10-17/62MNmnro.png")
$wsReceipts.Range("L5").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Sheet 3: "کاربران" (users) - refresh existing row 3 with the real user and
# append one new user row (4).
# ---------------------------------------------------------------------
$wsUsers = $wb.Worksheets.Item(3)

$wsUsers.Range("A3").Value = "Fzl2 M2"
# Phone numbers are text, not numbers - force text format so the digit
# string isn't auto-converted to a numeric cell, then drop back to the
# Normal cell style so no stray formatting is left behind.
$wsUsers.Range("B3").NumberFormat = "@"
$wsUsers.Range("B3").Value = "9130009997"
$wsUsers.Range("B3").Style = "Normal"
$wsUsers.Range("C3").Value = 120000
$wsUsers.Range("D3").Value = 12

$wsUsers.Range("A4").Value = "test test"
$wsUsers.Range("B4").NumberFormat = "@"
$wsUsers.Range("B4").Value = "9130009999"
$wsUsers.Range("B4").Style = "Normal"
$wsUsers.Range("C4").Value = 850000000
$wsUsers.Range("D4").Value = 123950
